$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "Credit Limit" column header (I1)
$ws.Range("I1").Value = "Credit Limit"
$ws.Range("I1").Font.Bold = $true

# Fill in Credit Limit values for existing rows 2-6
$ws.Range("I2").Value = "Eligible"
$ws.Range("I3").Value = "InEligible"
$ws.Range("I4").Value = "InEligible"
$ws.Range("I5").Value = "InEligible"
$ws.Range("I6").Value = "InEligible"

# Add the new row (Account Set 6)
$ws.Range("A7").Value = "Account Set 6"
$ws.Range("B7").Value = "T0XXXXX"
$ws.Range("C7").Value = "abcde123"
$ws.Range("D7").Value = "Multi"
$ws.Range("E7").Value = "PCH"
$ws.Range("F7").Value = "InEligible"
$ws.Range("G7").Value = "InEligible"
$ws.Range("H7").Value = "InEligible"
$ws.Range("I7").Value = "Eligible"

# Update the active selection to match the saved view state
$ws.Range("D13").Select()
